$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update confusion matrix values (row 3 = True Label 1)
$ws.Range("C3").Value = "16 (0.7619)"
$ws.Range("D3").Value = "5 (0.2381)"

# Row 4 values remain the same ("1 (0.0714)" and "13 (0.9286)") but
# the underlying shared-string indices still shuffle in the source
# workbook; re-asserting the same text keeps content correct either way.
$ws.Range("C4").Value = "1 (0.0714)"
$ws.Range("D4").Value = "13 (0.9286)"

# Update selection / active cell
$ws.Range("H3").Select()

# Update window scroll position (yWindow)
$excel.ActiveWindow.Top = 3600
